$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "TALAD1"
$ws.Range("B2").Value = "DEWALT"
$ws.Range("C2").Value = "TALADRO1"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "ricky"

# Row 3
$ws.Range("A3").Value = "CAD1"
$ws.Range("B3").Value = "Berger"
$ws.Range("C3").Value = "CADENA1"
$ws.Range("D3").Value = 150
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "ricky"

# Row 4
$ws.Range("A4").Value = "PINT2"
$ws.Range("B4").Value = "SICA"
$ws.Range("C4").Value = "PINTURA2"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "ricky"

# Row 5 (new)
$ws.Range("A5").Value = "AMOLA1"
$ws.Range("B5").Value = "DEWALT"
$ws.Range("C5").Value = "AMOLAD1"
$ws.Range("D5").Value = 100
$ws.Range("F5").Value = "ricky"

# Row 6 (new)
$ws.Range("A6").Value = "FOC1"
$ws.Range("B6").Value = "SICA"
$ws.Range("C6").Value = "FOCO1"
$ws.Range("D6").Value = 0
$ws.Range("F6").Value = "ferreteria_general"

# Row 7 (new) - keep an empty (but present) cell in column E, matching source row
$ws.Range("A7").Value = "AMOLA2"
$ws.Range("B7").Value = "DEWALT"
$ws.Range("C7").Value = "AMOLA2"
$ws.Range("D7").Value = 0
$ws.Range("E7").Font.Bold = $false
$ws.Range("F7").Value = "ricky"
